$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(239, 44313, 1, 6, 99.81700216270171),
    @(240, 44314, 1, 6, 99.81700216270171),
    @(241, 44315, 3, 9, 149.7255032440526),
    @(242, 44316, 2, 10, 166.3616702711695),
    @(243, 44317, 3, 13, 216.2701713525204),
    @(244, 44318, 0, 11, 182.9978372982865)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Copy style from the row above (row 238, which has the correct date style) onto the new rows A column
$ws.Range("A238").Copy()
$ws.Range("A239:A244").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
